# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. "58.925.04",
# "560.41") that must stay text, not get auto-converted to a number by
# Excel's type inference. Force the column to Text format first so the
# values below round-trip exactly as strings.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value2 = "58.925.04"
$ws.Range("E2").Value2 = "  +2.47%  "

# Row 3 - Ethereum
$ws.Range("D3").Value2 = "2.997.58"
$ws.Range("E3").Value2 = "  +1.68%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value2 = "  -0.10%  "

# Row 5 - BNB
$ws.Range("D5").Value2 = "560.41"
$ws.Range("E5").Value2 = "  +1.02%  "

# Row 6 - Solana
$ws.Range("D6").Value2 = "137.18"
$ws.Range("E6").Value2 = "  +4.05%  "

# Row 7 - USDC
$ws.Range("E7").Value2 = "  -0.23%  "

# Row 8 - XRP
$ws.Range("E8").Value2 = "  +1.48%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value2 = "2.985.21"
$ws.Range("E9").Value2 = "  +1.44%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value2 = "0.132"
$ws.Range("E10").Value2 = "  +3.85%  "

# Row 11 - Toncoin
$ws.Range("D11").Value2 = "5.17"
$ws.Range("E11").Value2 = "  +7.60%  "

# Row 12 - Cardano
$ws.Range("D12").Value2 = "0.456"
$ws.Range("E12").Value2 = "  +1.90%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value2 = "  +3.53%  "

# Row 14 - Avalanche
$ws.Range("D14").Value2 = "33.59"
$ws.Range("E14").Value2 = "  +2.62%  "

# Row 15 - TRON
$ws.Range("E15").Value2 = "  +2.21%  "

# Row 16 - Polkadot
$ws.Range("E16").Value2 = "  +7.44%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value2 = "3.484.81"
$ws.Range("E17").Value2 = "  +1.45%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value2 = "2.989.13"
$ws.Range("E18").Value2 = "  +1.51%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value2 = "58.935.10"
$ws.Range("E19").Value2 = "  +2.43%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value2 = "427.38"
$ws.Range("E20").Value2 = "  +2.71%  "

# Row 21 - Chainlink
$ws.Range("D21").Value2 = "13.67"
$ws.Range("E21").Value2 = "  +4.10%  "

# Row 22 - Polygon
$ws.Range("D22").Value2 = "0.724"
$ws.Range("E22").Value2 = "  +6.05%  "

# Row 23 - Uniswap
$ws.Range("D23").Value2 = "7.10"
$ws.Range("E23").Value2 = "  +2.13%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value2 = "13.27"
$ws.Range("E24").Value2 = "  +2.08%  "

# Row 25 - Litecoin
$ws.Range("D25").Value2 = "80.41"
$ws.Range("E25").Value2 = "  +1.55%  "

# Row 26 - Dai
$ws.Range("E26").Value2 = "  -0.02%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value2 = "  +0.08%  "

# Row 28 - ImmutableX
$ws.Range("D28").Value2 = "2.18"
$ws.Range("E28").Value2 = "  +10.33%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value2 = "2.54"
$ws.Range("E29").Value2 = "  +2.02%  "

# Row 30 - RenderToken
$ws.Range("D30").Value2 = "7.79"
$ws.Range("E30").Value2 = "  +3.63%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value2 = "25.71"
$ws.Range("E31").Value2 = "  +2.34%  "

# Row 32 - NEARProtocol
$ws.Range("D32").Value2 = "6.04"
$ws.Range("E32").Value2 = "  -0.13%  "

# Row 33 - Hedera
$ws.Range("D33").Value2 = "0.0993"
$ws.Range("E33").Value2 = "  -2.34%  "

# Row 34 - was Mantle, now Filecoin
$ws.Range("B34").Value2 = "Filecoin"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value2 = "5.94"
$ws.Range("E34").Value2 = "  +5.54%  "

# Row 35 - was Filecoin, now Mantle
$ws.Range("B35").Value2 = "Mantle"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value2 = "0.993"
$ws.Range("E35").Value2 = "  +5.48%  "

# Row 36 - PEPE
$ws.Range("D36").Value2 = "0.0₃0758"
$ws.Range("E36").Value2 = "  +10.71%  "

# Row 37 - Stacks
$ws.Range("D37").Value2 = "2.09"
$ws.Range("E37").Value2 = "  -1.26%  "

# Row 38 - OKB
$ws.Range("D38").Value2 = "48.78"
$ws.Range("E38").Value2 = "  +0.48%  "

# Row 39 - Cosmos
$ws.Range("E39").Value2 = "  +2.67%  "

# Row 40 - dogwifhat
$ws.Range("D40").Value2 = "2.72"
$ws.Range("E40").Value2 = "  +6.87%  "

# Row 41 - Bittensor
$ws.Range("D41").Value2 = "397.72"
$ws.Range("E41").Value2 = "  +5.07%  "

# Row 42 - VeChain
$ws.Range("E42").Value2 = "  +0.45%  "

# Row 43 - Maker
$ws.Range("D43").Value2 = "2.749.43"
$ws.Range("E43").Value2 = "  +3.41%  "

# Row 44 - Kaspa
$ws.Range("E44").Value2 = "  -0.87%  "

# Row 45 - TheGraph
$ws.Range("D45").Value2 = "0.251"
$ws.Range("E45").Value2 = "  +4.71%  "

# Row 46 - was USDe, now Arweave
$ws.Range("B46").Value2 = "Arweave"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value2 = "35.59"
$ws.Range("E46").Value2 = "  +25.81%  "

# Row 47 - was Arweave, now USDe
$ws.Range("B47").Value2 = "USDe"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value2 = "0.999"
$ws.Range("E47").Value2 = "  -0.04%  "

# Row 48 - Monero
$ws.Range("D48").Value2 = "123.05"
$ws.Range("E48").Value2 = "  +0.17%  "

# Row 49 - Stellar
$ws.Range("E49").Value2 = "  +0.93%  "

# Row 50 - Fetch.AI
$ws.Range("D50").Value2 = "2.00"
$ws.Range("E50").Value2 = "  +0.65%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").Value2 = "23.33"
$ws.Range("E51").Value2 = "  -0.19%  "
